$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HBNC")

# Insert two new columns before column D (shifts existing D:K -> F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Helper to set a pair of values for a given row into the new D and E columns
function Set-DE($row, $dVal, $eVal) {
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal
}

# Date header rows
Set-DE 7 43465 43373
Set-DE 38 43465 43373
Set-DE 80 43465 43373

# Income statement block
Set-DE 8 43700 42300
Set-DE 17 10400 9700
Set-DE 18 33300 32600
Set-DE 20 -17600 -17000
Set-DE 21 17400 17400
Set-DE 23 15700 15700
Set-DE 24 2500 2600
Set-DE 26 13100 13100
Set-DE 27 13100 13100
Set-DE 29 0 "NA"
Set-DE 32 17600 17000
Set-DE 33 13100 13100
Set-DE 35 13100 13100

# Balance sheet block
Set-DE 41 58500 69700
Set-DE 42 34000 23200
Set-DE 48 74300 75300
Set-DE 49 130300 130800
Set-DE 54 4246700 4150600
Set-DE 57 2000 1700
Set-DE 61 588200 515500
Set-DE 66 3754700 3673000
Set-DE 72 224000 214800
Set-DE 76 492000 477600

# Cash flow block
Set-DE 81 13100 13100
Set-DE 83 1800 1700
Set-DE 89 7800 19300
Set-DE 91 -100 -1400
Set-DE 94 -82400 -80100
Set-DE 96 -3900 -3800
Set-DE 100 79800 61500
Set-DE 102 5300 700

$wb.Save()
